# "Update countries & provincias Spain"
# This workbook's "Pais" sheet is a COVID-19 dashboard table (header row 3,
# data rows 4:219) sorted descending by column B ("Casos totales"). The
# source feed was re-pulled later the same day, so most rows below just get
# refreshed totals; a handful of countries with close case counts leap-frog
# their neighbour in the ranking, which is why those rows show both a new
# country name (col A) AND new stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer/title cell with the "last updated" timestamp.
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 01:29"

# --- Plain refreshes (ranking order unchanged) ---------------------------

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 2355664
$ws.Cells.Item(4, 3).Value = 25086
$ws.Cells.Item(4, 4).Value = 978399
$ws.Cells.Item(4, 5).Value = 1255022
$ws.Cells.Item(4, 7).Value = 263
$ws.Cells.Item(4, 8).Value = 122243

# Row 5 - Brasil
$ws.Cells.Item(5, 2).Value = 1086990
$ws.Cells.Item(5, 3).Value = 16851
$ws.Cells.Item(5, 5).Value = 493145
$ws.Cells.Item(5, 7).Value = 601
$ws.Cells.Item(5, 8).Value = 50659

# Row 25 - Colombia
$ws.Cells.Item(25, 2).Value = 68652
$ws.Cells.Item(25, 3).Value = 3019
$ws.Cells.Item(25, 4).Value = 27360
$ws.Cells.Item(25, 5).Value = 39055
$ws.Cells.Item(25, 7).Value = 111
$ws.Cells.Item(25, 8).Value = 2237

# Row 53 - Nigeria
$ws.Cells.Item(53, 2).Value = 20244
$ws.Cells.Item(53, 3).Value = 436
$ws.Cells.Item(53, 4).Value = 6879
$ws.Cells.Item(53, 5).Value = 12847
$ws.Cells.Item(53, 7).Value = 12
$ws.Cells.Item(53, 8).Value = 518

# Row 54 - Japon
$ws.Cells.Item(54, 2).Value = 17864
$ws.Cells.Item(54, 3).Value = 65
$ws.Cells.Item(54, 4).Value = 16108
$ws.Cells.Item(54, 5).Value = 803
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 953

# Row 67 - Chequia
$ws.Cells.Item(67, 2).Value = 10498
$ws.Cells.Item(67, 3).Value = 50
$ws.Cells.Item(67, 4).Value = 7499
$ws.Cells.Item(67, 5).Value = 2663
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 336

# Row 70 - Noruega
$ws.Cells.Item(70, 2).Value = 8745
$ws.Cells.Item(70, 3).Value = 3
$ws.Cells.Item(70, 5).Value = 363

# Row 131 - Burkina Faso
$ws.Cells.Item(131, 2).Value = 903
$ws.Cells.Item(131, 3).Value = 2
$ws.Cells.Item(131, 5).Value = 36

# Row 134 - Uruguay
$ws.Cells.Item(134, 2).Value = 876
$ws.Cells.Item(134, 3).Value = 17
$ws.Cells.Item(134, 4).Value = 814
$ws.Cells.Item(134, 5).Value = 37

# --- Ranking swaps (country names exchange rows, each row also gets the ---
# --- refreshed stats belonging to its new occupant) -----------------------

# Rows 65/66: Camerun overtakes Argelia
$ws.Cells.Item(65, 1).Value = "Camerun"
$ws.Cells.Item(65, 2).Value = 11892
$ws.Cells.Item(65, 3).Value = 282
$ws.Cells.Item(65, 4).Value = 7710
$ws.Cells.Item(65, 5).Value = 3879
$ws.Cells.Item(65, 7).Value = 2
$ws.Cells.Item(65, 8).Value = 303

$ws.Cells.Item(66, 1).Value = "Argelia"
$ws.Cells.Item(66, 2).Value = 11771
$ws.Cells.Item(66, 3).Value = 140
$ws.Cells.Item(66, 4).Value = 8422
$ws.Cells.Item(66, 5).Value = 2504
$ws.Cells.Item(66, 7).Value = 8
$ws.Cells.Item(66, 8).Value = 845

# Rows 150/151: Libia overtakes Togo
$ws.Cells.Item(150, 1).Value = "Libia"
$ws.Cells.Item(150, 2).Value = 571
$ws.Cells.Item(150, 3).Value = 27
$ws.Cells.Item(150, 4).Value = 103
$ws.Cells.Item(150, 5).Value = 458
$ws.Cells.Item(150, 8).Value = 10

$ws.Cells.Item(151, 1).Value = "Togo"
$ws.Cells.Item(151, 2).Value = 569
$ws.Cells.Item(151, 3).Value = 8
$ws.Cells.Item(151, 4).Value = 375
$ws.Cells.Item(151, 5).Value = 181
$ws.Cells.Item(151, 8).Value = 13

# Rows 202/203: Fiyi/Dominica tie-break swap (stats identical, name only)
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(203, 1).Value = "Dominica"

# Rows 207/208: Groenlandia/Islas Malvinas tie-break swap (stats identical)
$ws.Cells.Item(207, 1).Value = "Groenlandia"
$ws.Cells.Item(208, 1).Value = "Islas Malvinas"

# Rows 213/214: Papua Nueva Guinea overtakes Islas Virgenes Britanicas
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 8).Value = 1
